$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.05115733333333
$ws.Range("H2").Value = 126.153472
$ws.Range("I2").Value = 0.1594435451835853
$ws.Range("J2").Value = 0.1594435451835853
$ws.Range("O2").Value = 0.7426786721750401
$ws.Range("P2").Value = 0.7426786721750401
$ws.Range("Q2").Value = 26.99504882528711
$ws.Range("R2").Value = 242.955439427584
$ws.Range("S2").Value = 0.1184153204238261
$ws.Range("T2").Value = 0.1184153204238261
$ws.Range("G3").Value = 42.05115733333333
$ws.Range("H3").Value = 126.153472
$ws.Range("I3").Value = 0.1594435451835853
$ws.Range("J3").Value = 0.1594435451835853
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.2224236666666667
$ws.Range("N3").Value = 0.6672709999999999
$ws.Range("O3").Value = 0.2573213278249599
$ws.Range("P3").Value = 0.2573213278249599
$ws.Range("Q3").Value = 9.353172601656887
$ws.Range("R3").Value = 84.17855341491199
$ws.Range("S3").Value = 0.04102822475975915
$ws.Range("T3").Value = 0.04102822475975915
$ws.Range("G4").Value = 57.66057933333332
$ws.Range("I4").Value = 0.2186291119973147
$ws.Range("J4").Value = 0.2186291119973148
$ws.Range("O4").Value = 0.7426786721750401
$ws.Range("P4").Value = 0.7426786721750401
$ws.Range("Q4").Value = 37.01563174728177
$ws.Range("R4").Value = 333.1406857255359
$ws.Range("S4").Value = 0.1623711785969739
$ws.Range("T4").Value = 0.1623711785969739
$ws.Range("G5").Value = 57.66057933333332
$ws.Range("I5").Value = 0.2186291119973147
$ws.Range("J5").Value = 0.2186291119973148
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.2224236666666667
$ws.Range("N5").Value = 0.6672709999999999
$ws.Range("O5").Value = 0.2573213278249599
$ws.Range("P5").Value = 0.2573213278249599
$ws.Range("Q5").Value = 12.82507747744422
$ws.Range("R5").Value = 115.425697296998
$ws.Range("S5").Value = 0.0562579334003409
$ws.Range("T5").Value = 0.0562579334003409
$ws.Range("G6").Value = 99.15200299999999
$ws.Range("H6").Value = 297.456009
$ws.Range("I6").Value = 0.3759503393701321
$ws.Range("J6").Value = 0.3759503393701321
$ws.Range("O6").Value = 0.7426786721750401
$ws.Range("P6").Value = 0.7426786721750401
$ws.Range("Q6").Value = 63.65135544053867
$ws.Range("R6").Value = 572.8621989648479
$ws.Range("S6").Value = 0.2792102988471654
$ws.Range("T6").Value = 0.2792102988471654
$ws.Range("G7").Value = 99.15200299999999
$ws.Range("H7").Value = 297.456009
$ws.Range("I7").Value = 0.3759503393701321
$ws.Range("J7").Value = 0.3759503393701321
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.2224236666666667
$ws.Range("N7").Value = 0.6672709999999999
$ws.Range("O7").Value = 0.2573213278249599
$ws.Range("P7").Value = 0.2573213278249599
$ws.Range("Q7").Value = 22.05375206460433
$ws.Range("R7").Value = 198.483768581439
$ws.Range("S7").Value = 0.09674004052296667
$ws.Range("T7").Value = 0.09674004052296668
$ws.Range("G8").Value = 64.87322933333333
$ws.Range("H8").Value = 194.619688
$ws.Range("I8").Value = 0.2459770034489679
$ws.Range("J8").Value = 0.2459770034489679
$ws.Range("O8").Value = 0.7426786721750401
$ws.Range("P8").Value = 0.7426786721750401
$ws.Range("Q8").Value = 41.64584530754844
$ws.Range("R8").Value = 374.812607767936
$ws.Range("S8").Value = 0.1826818743070747
$ws.Range("T8").Value = 0.1826818743070747
$ws.Range("G9").Value = 64.87322933333333
$ws.Range("H9").Value = 194.619688
$ws.Range("I9").Value = 0.2459770034489679
$ws.Range("J9").Value = 0.2459770034489679
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.2224236666666667
$ws.Range("N9").Value = 0.6672709999999999
$ws.Range("O9").Value = 0.2573213278249599
$ws.Range("P9").Value = 0.2573213278249599
$ws.Range("Q9").Value = 14.42934153682755
$ws.Range("R9").Value = 129.864073831448
$ws.Range("S9").Value = 0.06329512914189314
$ws.Range("T9").Value = 0.06329512914189316
